# "Generate Report for Handback"
# Updates the localization-status report after a successful handback:
#  - Status moves from "Ready for handoff" to "Handed back: in sync with en-US"
#  - Latest Handback DateTime stamps are refreshed
#  - The stale "handback file is not the latest" Error Detail is cleared
#  - The Status / Error Detail columns are widened / narrowed to fit the new text

$wb = $excel.ActiveWorkbook

$newStatus = "Handed back: in sync with en-US"

# ---- Overview sheet -------------------------------------------------
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("E2").Value = $newStatus
$overview.Range("F2").Value = $newStatus
$overview.Range("E1").ColumnWidth = 29.9777047293527
$overview.Range("F1").ColumnWidth = 29.9777047293527

# ---- zh-cn sheet ------------------------------------------------------
$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("C2").Value = $newStatus
$zhcn.Range("K2").Value = "2016-08-22 22:48:28"
$zhcn.Range("P2").Value = ""
$zhcn.Range("C1").ColumnWidth = 29.9777047293527
$zhcn.Range("P1").ColumnWidth = 13.7470528738839

# ---- de-de sheet ------------------------------------------------------
$dede = $wb.Worksheets.Item("de-de")
$dede.Range("C2").Value = $newStatus
$dede.Range("K2").Value = "2016-08-22 22:48:35"
$dede.Range("P2").Value = ""
$dede.Range("C1").ColumnWidth = 29.9777047293527
$dede.Range("P1").ColumnWidth = 13.7470528738839
